$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated monthly stats: each row now reflects the following month's data,
# with a brand-new final row for 2024 October (rolling window update).
$data = @(
    @("2023 November", 55, 4),
    @("2023 December", 53, 7),
    @("2024 January", 54, 7),
    @("2024 February", 35, 13),
    @("2024 March", 36, 7),
    @("2024 April", 54, 13),
    @("2024 May", 38, 23),
    @("2024 June", 46, 25),
    @("2024 July", 34, 19),
    @("2024 August", 33, 23),
    @("2024 September", 24, 24),
    @("2024 October", 1, 8)
)

$row = 2
foreach ($entry in $data) {
    # Temporarily format column A as text so month labels like
    # "2023 November" aren't auto-converted into date serial values, then
    # clear the formatting back off so no stray number-format style is left
    # on the cell (matches the original, unstyled data rows).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $entry[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
